$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added to the dataset. It belongs
# chronologically right before the existing row 117, so insert a new row
# there (this pushes the former row 117 and everything below it down by
# one, which is exactly what the target diff shows).
$ws.Rows("117:117").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A117").Value = 6
$ws.Range("B117").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C117").Value = "Metropolitana"
$ws.Range("D117").Value = 44810
$ws.Range("E117").Value = 13
$ws.Range("F117").Value = 100112029
$ws.Range("G117").Value = "Orégano"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 48
$ws.Range("K117").Value = 15000
$ws.Range("L117").Value = 16000
$ws.Range("M117").Value = 15458
$ws.Range("N117").Value = "$/docena de atados"
$ws.Range("O117").Value = "Región Metropolitana"
$ws.Range("P117").Value = 5153
$ws.Range("Q117").Value = 3
$ws.Range("R117").Value = "Hortaliza"
